# Module 05 Release Candidate
# Applies the authored changes to the presentation:
#  1. Slide 1 (title slide): subtitle authors list updated; copyright year 2022 -> 2023
#  2. Slide 27: title "Odds and Ends" -> "Odds and Ends You Should Know About"
#  3. Slide 32: code sample "=>" -> "=> " (trailing space) in the TextBox 5 shape
#  4. Delete the last slide (id 552, "Learning Goals for this Lesson (expanded)")
#     and its notes page, since it was removed from the deck for the release
#     candidate.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 1 - title slide subtitle + copyright line
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

# Subtitle shape: "Jonathan Bell, Adeel Bhutta, Mitch Wand" -> "Adeel Bhutta, Jan Vitek, Mitch Wand"
$subtitle = $slide1.Shapes.Item(2)
$subTr = $subtitle.TextFrame.TextRange

$fullSub = $subTr.Text
$idx = $fullSub.IndexOf("Jonathan Bell, ")
if ($idx -ge 0) {
    $delRange = $subTr.Characters($idx + 1, 15)
    $delRange.Text = ""
}

$fullSub2 = $subTr.Text
$idx2 = $fullSub2.IndexOf("Bhutta")
if ($idx2 -ge 0) {
    $bhutta = $subTr.Characters($idx2 + 1, 6)
    $bhutta.InsertAfter(", Jan Vitek")
}

# Re-set the "Vitek" and "Bhutta" substrings in place so they become their own
# runs (mirrors the authored run layout without disturbing the rest of the text).
$fullSub3 = $subTr.Text
$idxVitek = $fullSub3.IndexOf("Vitek")
if ($idxVitek -ge 0) {
    $vitekRange = $subTr.Characters($idxVitek + 1, 5)
    $vitekRange.Text = "Vitek"
}

$fullSub4 = $subTr.Text
$idxBhutta = $fullSub4.IndexOf("Bhutta")
if ($idxBhutta -ge 0) {
    $bhuttaRange = $subTr.Characters($idxBhutta + 1, 6)
    $bhuttaRange.Text = "Bhutta"
}

# Copyright rectangle: "(c) 2022 Released under the" -> "(c) 2023 Released under the"
# Built from char codes (c)=U+00A9, nbsp=U+00A0) so the non-breaking space already
# present in the run is never round-tripped through .Text (which is lossy here).
$copyrightShape = $slide1.Shapes.Item(4)
$copyrightTr = $copyrightShape.TextFrame.TextRange
$copyrightRun = $copyrightTr.Characters(1, 26)
$newCopyrightText = [char]0x00A9 + " 2023 Released under the" + [char]0x00A0
$copyrightRun.Text = $newCopyrightText

# ---------------------------------------------------------------------------
# 2. Slide 27 - title update
# ---------------------------------------------------------------------------
$slide27 = $p.Slides.Item(27)
$titleShape = $slide27.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Odds and Ends You Should Know About"

# ---------------------------------------------------------------------------
# 3. Slide 32 - code sample tweak ("=>" -> "=> ")
# ---------------------------------------------------------------------------
$slide32 = $p.Slides.Item(32)
$codeShape = $slide32.Shapes.Item(4)
$codeTr = $codeShape.TextFrame.TextRange
$fullCode = $codeTr.Text
$idxArrow = $fullCode.IndexOf("=>")
if ($idxArrow -ge 0) {
    $arrowRange = $codeTr.Characters($idxArrow + 1, 2)
    $arrowRange.Text = "=> "
}

# ---------------------------------------------------------------------------
# 4. Remove the final slide ("Learning Goals for this Lesson (expanded)")
# ---------------------------------------------------------------------------
$lastSlide = $p.Slides.Item($p.Slides.Count)
$lastSlide.Delete()
